$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 116, shifting existing rows 116..215 down to 117..216
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row 116 with the new record
$ws.Cells.Item(116, 1).Value = 7
$ws.Cells.Item(116, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(116, 3).Value = "Ñuble"
$ws.Cells.Item(116, 4).Value = 44651
$ws.Cells.Item(116, 5).Value = 16
$ws.Cells.Item(116, 6).Value = 100112043
$ws.Cells.Item(116, 7).Value = "Pepino ensalada"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 60
$ws.Cells.Item(116, 11).Value = 16000
$ws.Cells.Item(116, 12).Value = 17000
$ws.Cells.Item(116, 13).Value = 16500
$ws.Cells.Item(116, 14).Value = "$/caja 80 unidades"
$ws.Cells.Item(116, 15).Value = "Región del Maule"
$ws.Cells.Item(116, 16).Value = 206
$ws.Cells.Item(116, 17).Value = 80
$ws.Cells.Item(116, 18).Value = "Hortaliza"
